# SummaryStats: the "last" sampling result recorded for each well must be a
# detected result (not a non-detect / estimated value), so the Last
# Cr / Last Date columns (X/Y) for several wells are updated to the most
# recent *detected* value and its date. Where no detected result exists at
# all, both cells are replaced with "No Detect Data".
#
# Every value in these columns is stored as literal text (not a true number
# or date), even though many of them look numeric/date-like. A plain
# `Range.Value = "..."` assignment would let Excel's normal type inference
# convert a string like "216.0" into the number 216, or "1998-06-09" into a
# date serial - which also reformats the cell. To avoid that, each target
# cell is temporarily switched to the Text number format before the literal
# is written, and then its original formatting is restored afterwards by
# copy/pasting the format from an untouched cell in the same column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Intermediate for Mapping")

# Untouched reference cells (same columns, unaffected row) used purely to
# restore the original cell style/format after forcing a text entry.
$xFormatSource = $ws.Range("X15")
$yFormatSource = $ws.Range("Y15")

function Set-TextValue($range, $value, $formatSource) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $formatSource.Copy()
    $range.PasteSpecial(-4122) # xlPasteFormats
}

Set-TextValue $ws.Range("X3") "216.0" $xFormatSource
Set-TextValue $ws.Range("Y3") "1998-06-09" $yFormatSource

Set-TextValue $ws.Range("X4") "1.3962" $xFormatSource
Set-TextValue $ws.Range("Y4") "2011-06-03" $yFormatSource

Set-TextValue $ws.Range("X8") "1.3" $xFormatSource

Set-TextValue $ws.Range("X9") "No Detect Data" $xFormatSource
Set-TextValue $ws.Range("Y9") "No Detect Data" $yFormatSource

Set-TextValue $ws.Range("X12") "0.418" $xFormatSource
Set-TextValue $ws.Range("Y12") "2012-09-06" $yFormatSource

Set-TextValue $ws.Range("X13") "1.9" $xFormatSource
Set-TextValue $ws.Range("Y13") "2005-03-24" $yFormatSource

Set-TextValue $ws.Range("X14") "6.01" $xFormatSource

$excel.CutCopyMode = 0
